$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: Excel's normal Range.Value assignment "smart-parses" text
# that looks numeric/date-like (phone numbers, DDD codes, ISO dates), which
# would corrupt these text-only columns. Force the destination cell to the
# "Text" number format before the write so the literal string is kept, then
# copy the (unformatted/General) look of an untouched sibling cell back on
# top so the stored style index is not altered.
function Set-TextValue($ws, $targetAddr, $likeAddr, $value) {
    $target = $ws.Range($targetAddr)
    $target.NumberFormat = "@"
    $target.Value = $value
    $ws.Range($likeAddr).Copy()
    $target.PasteSpecial(-4122)  # xlPasteFormats
}

# Several phone-number rows were removed from the list; the remaining
# entries shift up, so rows 4-6 now hold what used to be rows 7, 8 and 11.
Set-TextValue $ws "A4" "A2" "+5511976440031"
Set-TextValue $ws "B4" "B2" "11"
Set-TextValue $ws "C4" "C2" "2024-09-16"

Set-TextValue $ws "A5" "A2" "+556293286544"
Set-TextValue $ws "B5" "B2" "62"
Set-TextValue $ws "C5" "C2" "2024-09-16"

Set-TextValue $ws "A6" "A2" "+5521981400589"
Set-TextValue $ws "B6" "B2" "21"
Set-TextValue $ws "C6" "C2" "2024-09-09"

# Remove the now-obsolete trailing rows (old rows 7-11), shifting cells
# below them up so the used range shrinks down to A1:C6.
$ws.Range("A7:C11").Delete(-4162)  # xlShiftUp
